$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "28.216.93"
Set-TextValue "E2" "  +0.94%  "
Set-TextValue "D3" "1.869.76"
Set-TextValue "E3" "  +3.65%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.55%  "
Set-TextValue "D5" "311.44"
Set-TextValue "E5" "  +1.90%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  -0.73%  "
Set-TextValue "D7" "0.5066"
Set-TextValue "E7" "  +2.09%  "
Set-TextValue "D8" "0.3926"
Set-TextValue "E8" "  +1.86%  "
Set-TextValue "D9" "0.09675"
Set-TextValue "E9" "  +3.67%  "
Set-TextValue "D10" "1.141"
Set-TextValue "E10" "  +4.17%  "
Set-TextValue "D11" "40.91"
Set-TextValue "E11" "  +0.76%  "
Set-TextValue "D12" "6.514"
Set-TextValue "E12" "  +3.04%  "
Set-TextValue "E13" "  +2.01%  "
Set-TextValue "D14" "1.874.48"
Set-TextValue "E14" "  +3.49%  "
Set-TextValue "D15" "7.440"
Set-TextValue "E15" "  +3.05%  "
Set-TextValue "E16" "  -0.52%  "
Set-TextValue "D17" "0.00001130"
Set-TextValue "E17" "  +1.44%  "
Set-TextValue "D18" "93.00"
Set-TextValue "E18" "  +1.18%  "
Set-TextValue "D19" "0.06584"
Set-TextValue "D20" "17.56"
Set-TextValue "E20" "  +2.84%  "
Set-TextValue "E21" "  -0.71%  "
Set-TextValue "D22" "6.166"
Set-TextValue "E22" "  +3.53%  "
Set-TextValue "D23" "28.279.87"
Set-TextValue "E23" "  +0.98%  "
Set-TextValue "D24" "11.38"
Set-TextValue "E24" "  +3.41%  "
Set-TextValue "D25" "2.289"
Set-TextValue "E25" "  +3.23%  "
Set-TextValue "D26" "2.552"
Set-TextValue "E26" "  +8.16%  "
Set-TextValue "D27" "2.088.33"
Set-TextValue "E27" "  +3.50%  "
Set-TextValue "D28" "21.22"
Set-TextValue "E28" "  +3.87%  "
Set-TextValue "D29" "158.51"
Set-TextValue "E29" "  +0.60%  "
Set-TextValue "E30" "  +0.86%  "
Set-TextValue "E31" "  -1.02%  "
Set-TextValue "E32" "  +2.62%  "
Set-TextValue "D33" "5.639"
Set-TextValue "E33" "  +1.81%  "
Set-TextValue "D34" "3.623"
Set-TextValue "E34" "  -0.15%  "
Set-TextValue "D35" "9.549"
Set-TextValue "E35" "  +7.17%  "
Set-TextValue "D36" "0.06723"
Set-TextValue "E36" "  -1.60%  "
Set-TextValue "D37" "0.02387"
Set-TextValue "E37" "  +3.60%  "
Set-TextValue "D38" "0.2190"
Set-TextValue "E38" "  +2.91%  "
Set-TextValue "D41" "4.979"
Set-TextValue "E41" "  +1.07%  "
Set-TextValue "D42" "1.185"
Set-TextValue "E42" "  +3.19%  "
Set-TextValue "D43" "1.000"
Set-TextValue "E43" "  -0.72%  "
Set-TextValue "D44" "13.59"
Set-TextValue "E44" "  +3.97%  "
Set-TextValue "D45" "0.6025"
Set-TextValue "E45" "  +2.53%  "
Set-TextValue "D46" "3.660"
Set-TextValue "E46" "  -0.24%  "
Set-TextValue "D47" "1.260"
Set-TextValue "E47" "  -2.16%  "
Set-TextValue "E48" "  +3.10%  "
Set-TextValue "D49" "124.07"
Set-TextValue "E49" "  +0.20%  "
Set-TextValue "D50" "1.197"
Set-TextValue "E50" "  +2.54%  "
Set-TextValue "D51" "0.06846"
Set-TextValue "E51" "  +1.71%  "

# Row 39 and 40 swap (Aptos <-> TheSandbox) with updated values
Set-TextValue "B39" "TheSandbox"
Set-TextValue "C39" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D39" "0.6378"
Set-TextValue "E39" "  +4.14%  "

Set-TextValue "B40" "Aptos"
Set-TextValue "C40" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D40" "11.51"
Set-TextValue "E40" "  +1.42%  "